$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 2.35

# Row 4
$ws.Range("G4").Value = 1.87
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.63
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Y4").Value = 1.5
$ws.Range("Z4").Value = 2.5
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 8.5
$ws.Range("AF4").Value = 17
$ws.Range("AJ4").Value = 6.5
$ws.Range("AK4").Value = 17
$ws.Range("AN4").Value = 9.5
$ws.Range("AO4").Value = 19
$ws.Range("AP4").Value = 13
$ws.Range("AQ4").Value = 41
$ws.Range("AR4").Value = 34

# Row 5
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5.25
$ws.Range("L5").Value = 6
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("AA5").Value = 2.1
$ws.Range("AB5").Value = 1.67
$ws.Range("AH5").Value = 34
